# Add 2022-Q3 data:
#  1. Insert a new "2022-Q3" worksheet right after "总计", populated with the
#     new quarter's fund-holding detail rows.
#  2. Insert a new summary row into "总计" for 2022-Q3, shifting the existing
#     quarters down and renumbering the index column.
# The other quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3) need no edits
# of their own - they simply shift one tab position to the right because the
# new sheet is inserted before them.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) New "2022-Q3" sheet, positioned right after "总计"
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ3.Name = "2022-Q3"

$headerCols = @("B","C","D","E","F","G","H")
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $headerCols[$i]
    $cell = $wsQ3.Range("$col" + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# code, name, size, total position, position ratio, held value(billion), rank
$data = @(
    @("012368", "摩根士丹利华鑫优享臻选六个月持有期混合A", "4.42", "94.05", "5.73", "0.2533", 8),
    @("000309", "大摩品质生活精选股票",                      "3.42", "94.27", "6.01", "0.2055", 8),
    @("233006", "大摩领先优势混合",                          "3.47", "94.24", "5.81", "0.2016", 6),
    @("010322", "大摩新兴产业股票",                          "1.98", "94.19", "5.68", "0.1125", 8),
    @("012369", "摩根士丹利华鑫优享臻选六个月持有期混合C", "0.30", "94.05", "5.73", "0.0172", 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $acell = $wsQ3.Range("A$r")
    $acell.Value = $i
    $acell.Font.Bold = $true
    $acell.HorizontalAlignment = -4108
    $acell.VerticalAlignment = -4160

    $bcell = $wsQ3.Range("B$r")
    $bcell.NumberFormat = "@"
    $bcell.Value = $row[0]

    $wsQ3.Range("C$r").Value = $row[1]

    $dcell = $wsQ3.Range("D$r")
    $dcell.NumberFormat = "@"
    $dcell.Value = $row[2]

    $ecell = $wsQ3.Range("E$r")
    $ecell.NumberFormat = "@"
    $ecell.Value = $row[3]

    $fcell = $wsQ3.Range("F$r")
    $fcell.NumberFormat = "@"
    $fcell.Value = $row[4]

    $gcell = $wsQ3.Range("G$r")
    $gcell.NumberFormat = "@"
    $gcell.Value = $row[5]

    $wsQ3.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) "总计" summary sheet: insert the 2022-Q3 row, shift the rest down
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$a2 = $wsTotal.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.79

# renumber the index column for the rows that shifted down
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
